$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("obj_output")
$ws2 = $wb.Worksheets.Item("obj_report")
$ws3 = $wb.Worksheets.Item("rel_report__output")

# ---------------------------------------------------------------------
# Introduce the three brand-new output names in the exact order that
# makes them land at shared-string indices 14, 15, 16 respectively:
#   14 = node_injection
#   15 = node_slack_pos
#   16 = node_slack_neg
# We use the cells that will end up holding them first (sheet1 B5,
# B4, B11) purely to control shared-string allocation order; their
# final values are (re)written below in normal sheet order as well.
# ---------------------------------------------------------------------
$ws1.Range("B5").Value = "node_injection"
$ws1.Range("B4").Value = "'node_slack_pos"
$ws1.Range("B11").Value = "node_slack_neg"

# ===========================================================================
# Sheet "obj_output" (sheet1): reorder / extend the list of output names
# ===========================================================================
$ws1.Range("B3").Value = "units_started_up"
$ws1.Range("B4").Value = "'node_slack_pos"
$ws1.Range("B4").NumberFormat = "d-mmm"
$ws1.Range("B5").Value = "node_injection"
$ws1.Range("B8").Value = "'units_shut_down"
$ws1.Range("B8").NumberFormat = "d-mmm"

$ws1.Range("A9").Value = "output"
$ws1.Range("B9").Value = "node_state"

$ws1.Range("A10").Value = "output"
$ws1.Range("B10").Value = "units_on"

$ws1.Range("A11").Value = "output"
$ws1.Range("B11").Value = "node_slack_neg"

$null = $ws1.Activate()
$null = $ws1.Range("A9:A11").Select()

# ===========================================================================
# Sheet "obj_report" (sheet2): no data changes, it simply stops being the
# active/selected tab (handled implicitly below by activating sheet3 last).
# ===========================================================================

# ===========================================================================
# Sheet "rel_report__output" (sheet3): reorder / extend the list of outputs
# ===========================================================================
$ws3.Range("C3").Value = "units_started_up"
$ws3.Range("C4").Value = "node_slack_pos"
$ws3.Range("C5").Value = "'node_injection"
$ws3.Range("C5").NumberFormat = "d-mmm"
$ws3.Range("C6").Value = "units_available"
$ws3.Range("C7").Value = "connection_flow"
$ws3.Range("C8").Value = "'units_shut_down"
$ws3.Range("C8").NumberFormat = "d-mmm"

$ws3.Range("A9").Value = "report__output"
$ws3.Range("B9").Value = "'result_temp"
$ws3.Range("B9").NumberFormat = "d-mmm"
$ws3.Range("C9").Value = "node_state"

$ws3.Range("A10").Value = "report__output"
$ws3.Range("B10").Value = "'result_temp"
$ws3.Range("B10").NumberFormat = "d-mmm"
$ws3.Range("C10").Value = "units_on"

$ws3.Range("A11").Value = "report__output"
$ws3.Range("B11").Value = "'result_temp"
$ws3.Range("B11").NumberFormat = "d-mmm"
$ws3.Range("C11").Value = "node_slack_neg"

# Slightly narrow/widen the bestFit columns A and B to track the (very
# marginally) different text metrics after the edit.
$ws3.Columns.Item(1).ColumnWidth = 13
$ws3.Columns.Item(2).ColumnWidth = 10

$null = $ws3.Activate()
$null = $ws3.Range("E6").Select()
